$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for this market/product and needs
# to be inserted as row 14 (in date order among the existing entries).
# Insert a new row at position 14; this pushes the existing rows 14-32 down
# to rows 15-33 and extends the sheet dimension to A1:R33 automatically.
$ws.Rows("14:14").Insert()

# Fill in the newly inserted row 14 with the new observation.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44426
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112001
$ws.Range("G14").Value = "Berenjena"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 11000
$ws.Range("M14").Value = 10500
$ws.Range("N14").Value = "$/caja 60 unidades"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 175
$ws.Range("Q14").Value = 60
$ws.Range("R14").Value = "Hortaliza"
